# Validations & import fixes
# Applies the changes captured in the commit's OOXML diff:
#   - Product:    add ID 5 (no name yet) and a stray "LM" name (no ID yet)
#   - Module:     add "Wilin' Out" (ProductID 5) and a bad/unfinished row
#                 where the Name column was typed in as a number (1)
#   - System Area: clear two stray cells (B4, A6) and add three new rows
#                 (Claims / Detentions / Cost Allocation) under ModuleID 1
#   - Key Action: fix the casing typo "Create Customer Order" -> "create
#                 Customer Order" and add a new "Choo Choo!" row
#   - Input Parameter: no data changes
#   - Product becomes the active/selected sheet (was Input Parameter)

$wb = $excel.ActiveWorkbook

$wsProduct = $wb.Worksheets.Item("Product")
$wsModule = $wb.Worksheets.Item("Module")
$wsSystemArea = $wb.Worksheets.Item("System Area")
$wsKeyAction = $wb.Worksheets.Item("Key Action")
$wsInputParameter = $wb.Worksheets.Item("Input Parameter")

# --- Product -------------------------------------------------------------
$wsProduct.Range("A6").Value = 5
$wsProduct.Range("B7").Value = "LM"

# --- Module ----------------------------------------------------------------
$wsModule.Range("A8").Value = 7
$wsModule.Range("B8").Value = 5
$wsModule.Range("C8").Value = "Wilin' Out"

$wsModule.Range("A9").Value = 8
$wsModule.Range("B9").Value = 2
$wsModule.Range("C9").Value = 1

# --- System Area -----------------------------------------------------------
$wsSystemArea.Range("B4").ClearContents()
$wsSystemArea.Range("A6").ClearContents()

$wsSystemArea.Range("A8").Value = 7
$wsSystemArea.Range("B8").Value = 1
$wsSystemArea.Range("C8").Value = "Claims"

$wsSystemArea.Range("A9").Value = 8
$wsSystemArea.Range("B9").Value = 1
$wsSystemArea.Range("C9").Value = "Detentions"

$wsSystemArea.Range("A10").Value = 9
$wsSystemArea.Range("B10").Value = 1
$wsSystemArea.Range("C10").Value = "Cost Allocation"

# --- Key Action --------------------------------------------------------
$wsKeyAction.Range("C7").Value = "create Customer Order"

$wsKeyAction.Range("C8").Value = "Choo Choo!"
$wsKeyAction.Range("D8").Value = "Riding the train!"
$wsKeyAction.Range("E8").Value = 1

# --- Selections --------------------------------------------------------
# Set the non-active sheets' selection first so the very last Select()
# call (on Product) is the one that ends up driving the active tab.
$wsModule.Range("A10").Select()
$wsSystemArea.Range("A7").Select()
$wsKeyAction.Range("D18").Select()
$wsInputParameter.Range("F23").Select()

# Product becomes the selected/active sheet.
$wsProduct.Range("B8").Select()
